# Completa ratoli, afegir Community detection i Cliques. Apuntar tot aixo al excel
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Full1")

# --- Column H = "Community detection" ---------------------------------
$ws.Range("H3").Value  = "SI"
$ws.Range("H4").Value  = "SI"
$ws.Range("H5").Value  = "SI"
$ws.Range("H6").Value  = "SI"
$ws.Range("H7").Value  = "Tarda"
$ws.Range("H8").Value  = "Tarda"
$ws.Range("H9").Value  = "SI"
$ws.Range("H10").Value = "SI"
$ws.Range("H11").Value = "SI"
$ws.Range("H13").Value = "SI"
$ws.Range("H14").Value = "SI"

# --- Column J = "Cliques and cavities" ---------------------------------
$ws.Range("J3").Value  = "SI"
$ws.Range("J4").Value  = "SI"
$ws.Range("J5").Value  = "SI"
$ws.Range("J6").Value  = "SI"
$ws.Range("J7").Value  = "Tarda"
$ws.Range("J8").Value  = "Tarda"
$ws.Range("J9").Value  = "no"
$ws.Range("J10").Value = "SI"
$ws.Range("J11").Value = "SI"
$ws.Range("J12").Value = "SI"
$ws.Range("J13").Value = "SI"
$ws.Range("J14").Value = "SI"

# --- Mouse row (row 19) now filled in -----------------------------------
$ws.Range("C19").Value = "SI"
$ws.Range("D19").Value = "SI"
$ws.Range("E19").Value = "SI"
$ws.Range("F19").Value = "SI"
$ws.Range("G19").Value = "SI"
$ws.Range("I19").Value = "SI"

# --- Fly Opt - sw row correction (row 13) --------------------------------
$ws.Range("D13").Value = "NO"

# --- Selection moved to J2 and scrolled back to the top -----------------
[void]$ws.Range("J2").Select()
